# Weekly CompStat update: new crime data collected for week 7/14/2025 - 7/20/2025.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report number + date range) ---
# A8 holds "Volume 32   Number  28" as rich text runs; only the trailing
# issue number changes (28 -> 29). All runs share identical formatting
# (10pt Andale WT), so a plain replacement is visually/semantically identical.
$ws.Range("A8").Value = "Volume 32   Number  29"

# C9 holds "Report Covering the Week  7/7/2025  Through  7/13/2025"; the
# week has rolled forward by one week on each side.
$ws.Range("C9").Value = "Report Covering the Week  7/14/2025  Through  7/20/2025"

# --- Body table updates: new weekly complaint statistics ---
$ws.Range("G15").Value = 3
$ws.Range("L15").Value = -12.5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 20
$ws.Range("I16").Value = 47
$ws.Range("J16").Value = 28
$ws.Range("K16").Value = 67.857142857142
$ws.Range("L16").Value = -11.320754716981
$ws.Range("M16").Value = -74.731182795698
$ws.Range("N16").Value = -91.376146788990
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 13
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -35
$ws.Range("I17").Value = 113
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = -24.666666666666
$ws.Range("L17").Value = -18.115942028985
$ws.Range("M17").Value = -33.918128654970
$ws.Range("N17").Value = -45.145631067961
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 33.333333333333
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = 22.580645161290
$ws.Range("L18").Value = -1.298701298701
$ws.Range("M18").Value = -60.824742268041
$ws.Range("N18").Value = -91.402714932126
$ws.Range("G19").Value = 29
$ws.Range("H19").Value = -13.793103448275
$ws.Range("I19").Value = 197
$ws.Range("J19").Value = 204
$ws.Range("K19").Value = -3.431372549019
$ws.Range("L19").Value = -11.261261261261
$ws.Range("M19").Value = -20.883534136546
$ws.Range("N19").Value = -40.483383685800
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = -77.777777777777
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -23.529411764705
$ws.Range("I20").Value = 102
$ws.Range("J20").Value = 131
$ws.Range("K20").Value = -22.137404580152
$ws.Range("L20").Value = 18.604651162790
$ws.Range("M20").Value = -51.196172248803
$ws.Range("N20").Value = -94.432314410480
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = -39.130434782608
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -16.867469879518
$ws.Range("I21").Value = 543
$ws.Range("J21").Value = 588
$ws.Range("K21").Value = -7.653061224489
$ws.Range("L21").Value = -7.179487179487
$ws.Range("M21").Value = -47.230320699708
$ws.Range("N21").Value = -85.837245696400
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 28.571428571428
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 35.714285714285
$ws.Range("I24").Value = 320
$ws.Range("J24").Value = 382
$ws.Range("K24").Value = -16.230366492146
$ws.Range("L24").Value = -26.940639269406
$ws.Range("M24").Value = -31.182795698924
$ws.Range("C25").Value = 5
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 14
$ws.Range("H25").Value = 21.428571428571
$ws.Range("I25").Value = 85
$ws.Range("K25").Value = -10.526315789473
$ws.Range("L25").Value = 6.25
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 27
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = -41.304347826087
$ws.Range("I26").Value = 237
$ws.Range("J26").Value = 235
$ws.Range("K26").Value = 0.851063829787
$ws.Range("L26").Value = 26.063829787234
$ws.Range("M26").Value = -35.597826086956
$ws.Range("G27").Value = 3
$ws.Range("L27").Value = -41.666666666666
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 200
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = 100
$ws.Range("L31").Value = 300

# A handful of cells flip from a numeric value to the workbook's standard
# "no data" placeholders: the text "0" (shared string, style 13) for counts
# and "***.*" (shared string, style 13) for the associated % change, used
# throughout this sheet (see e.g. row 22) whenever the underlying count is
# zero/undefined for that period. Copy formatting+value from an existing
# placeholder cell so the style/shared-string plumbing matches exactly.
$zeroTemplate = $ws.Range("C22")
$starTemplate = $ws.Range("E22")

$zeroTemplate.Copy($ws.Range("D15"))
$starTemplate.Copy($ws.Range("E15"))

$zeroTemplate.Copy($ws.Range("D25"))
$starTemplate.Copy($ws.Range("E25"))

$zeroTemplate.Copy($ws.Range("D27"))
$starTemplate.Copy($ws.Range("E27"))

$zeroTemplate.Copy($ws.Range("C28"))
$zeroTemplate.Copy($ws.Range("D28"))
$starTemplate.Copy($ws.Range("E28"))
